$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row within the used columns (A:R) at row 377, shifting rows 377:461 down to 378:462
$ws.Range("A377:R377").Insert(-4121) # xlShiftDown

# Populate the new row 377 with the new data (matches old row 377 except for changed fields)
$ws.Cells.Item(377, 1).Value = 8
$ws.Cells.Item(377, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(377, 3).Value = "Coquimbo"
$ws.Cells.Item(377, 4).Value = 45211
$ws.Cells.Item(377, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(377, 5).Value = 4
$ws.Cells.Item(377, 6).Value = 100112031
$ws.Cells.Item(377, 7).Value = "Poroto verde"
$ws.Cells.Item(377, 8).Value = "Magnum"
$ws.Cells.Item(377, 9).Value = "Primera"
$ws.Cells.Item(377, 10).Value = 520
$ws.Cells.Item(377, 11).Value = 26000
$ws.Cells.Item(377, 12).Value = 27000
$ws.Cells.Item(377, 13).Value = 26500
$ws.Cells.Item(377, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(377, 15).Value = "Perú"
$ws.Cells.Item(377, 16).Value = 1060
$ws.Cells.Item(377, 17).Value = 25
$ws.Cells.Item(377, 18).Value = "Hortaliza"
